$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 8 and row 9 content for columns F..V (A..E - Indice/pais/torneio/temporada/data_partida - stay put) ---
$row8 = @()
$row9 = @()
for ($c = 6; $c -le 22; $c++) {
    $row8 += $ws.Cells.Item(8, $c).Value()
    $row9 += $ws.Cells.Item(9, $c).Value()
}
for ($i = 0; $i -lt $row8.Length; $i++) {
    $c = 6 + $i
    $ws.Cells.Item(8, $c).Value = $row9[$i]
    $ws.Cells.Item(9, $c).Value = $row8[$i]
}

# --- Swap row 22 and row 23 content for columns F..V ---
$row22 = @()
$row23 = @()
for ($c = 6; $c -le 22; $c++) {
    $row22 += $ws.Cells.Item(22, $c).Value()
    $row23 += $ws.Cells.Item(23, $c).Value()
}
for ($i = 0; $i -lt $row22.Length; $i++) {
    $c = 6 + $i
    $ws.Cells.Item(22, $c).Value = $row23[$i]
    $ws.Cells.Item(23, $c).Value = $row22[$i]
}

# --- Append new match row 46 ---
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "lebanon"
$ws.Cells.Item(46, 3).Value = "premier-league"
$ws.Cells.Item(46, 4).Value = "2023-2024"
$ws.Cells.Item(46, 5).Value = 45241.45833333334
$ws.Cells.Item(46, 6).Value = "Nejmeh SC"
$ws.Cells.Item(46, 7).Value = 2
$ws.Cells.Item(46, 8).Value = "Al Hikma"
$ws.Cells.Item(46, 9).Value = 1
$ws.Cells.Item(46, 10).Value = 1.16
$ws.Cells.Item(46, 11).Value = "09/11/2023 23:13"
$ws.Cells.Item(46, 12).Value = 1.2
$ws.Cells.Item(46, 13).Value = "11/11/2023 10:46"
$ws.Cells.Item(46, 14).Value = 6.19
$ws.Cells.Item(46, 15).Value = "09/11/2023 23:13"
$ws.Cells.Item(46, 16).Value = 6.19
$ws.Cells.Item(46, 17).Value = "11/11/2023 10:46"
$ws.Cells.Item(46, 18).Value = 12.58
$ws.Cells.Item(46, 19).Value = "09/11/2023 23:13"
$ws.Cells.Item(46, 20).Value = 13.15
$ws.Cells.Item(46, 21).Value = "11/11/2023 10:46"
$ws.Cells.Item(46, 22).Value = "https://www.betexplorer.com/football/lebanon/premier-league/nejmeh-sc-al-hikma/pKttIZCi/"

# Copy formatting for the new row's special-format columns (A: bordered/bold index style, E: datetime format)
# from the previous last row (45), matching the rest of the sheet.
$ws.Range("A45").Copy()
$ws.Range("A46").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E45").Copy()
$ws.Range("E46").PasteSpecial(-4122) # xlPasteFormats
